$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "33.808.27"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.779.45"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.27"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.01"
$ws.Range("E8").Value = "  -1.34%  "
$ws.Range("E9").Value = "  +1.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0675"
$ws.Range("E10").Value = "  -5.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0934"
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.036.42"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.19"
$ws.Range("E13").Value = "  +4.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.772.03"
$ws.Range("E14").Value = "  -1.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.840.16"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.610"
$ws.Range("E16").Value = "  -3.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.12"
$ws.Range("E17").Value = "  -2.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.55"
$ws.Range("E18").Value = "  -2.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "238.02"
$ws.Range("E19").Value = "  -3.50%  "
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.57"
$ws.Range("E22").Value = "  -2.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.00"
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("E24").Value = "  -2.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.16"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.07"
$ws.Range("E26").Value = "  -2.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.01"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("E31").Value = "  -2.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.58"
$ws.Range("E32").Value = "  -3.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.49"
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("E34").Value = "  -2.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.381.94"
$ws.Range("E35").Value = "  -2.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.636"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("E37").Value = "  -2.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0184"
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("B39").Value = "HuobiToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.39"
$ws.Range("E39").Value = "  +2.12%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.23"
$ws.Range("E40").Value = "  +4.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "78.41"
$ws.Range("E41").Value = "  -2.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.908"
$ws.Range("E42").Value = "  -3.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.55"
$ws.Range("E43").Value = "  +12.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.65"
$ws.Range("E44").Value = "  -2.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0142"
$ws.Range("E45").Value = "  +16.04%  "
$ws.Range("E46").Value = "  +1.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.08"
$ws.Range("E47").Value = "  +3.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "107.62"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.84"
$ws.Range("E49").Value = "  -2.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.937.68"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.06%  "
